$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.736.42"
$ws.Range("E2").Value = "  -1.95%  "

$ws.Range("D3").Value = "1.758.08"
$ws.Range("E3").Value = "  -2.32%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4533"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3733"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.35"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07558"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.128"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9994"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.203"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.329"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.52%  "

$ws.Range("D16").Value = "1.753.35"
$ws.Range("E16").Value = "  -2.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001074"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06213"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9989"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.187"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5313"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.26%  "

$ws.Range("D24").Value = "27.780.16"
$ws.Range("E24").Value = "  -1.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.314"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.366"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").Value = "1.952.31"
$ws.Range("E30").Value = "  -2.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "128.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.225"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09349"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.744"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.652"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2187"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02334"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.75%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6511"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.58%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06145"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.20%  "

$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.100"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.200"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.011"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.415"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9990"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.88%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.88%  "

$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.756"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.992"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06914"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.19%  "

